$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D/E columns stay text-formatted so numeric-looking values
# (e.g. "1.00", "169.20") are not coerced into numbers, matching the
# original inline-string cell content.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "36.648.72"
$ws.Range("E2").Value = "  +3.29%  "
$ws.Range("D3").Value = "2.066.11"
$ws.Range("E3").Value = "  +9.46%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "248.03"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("D6").Value = "0.666"
$ws.Range("E6").Value = "  -3.05%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "45.12"
$ws.Range("E8").Value = "  +5.19%  "
$ws.Range("D9").Value = "60.85"
$ws.Range("E9").Value = "  +7.60%  "
$ws.Range("D10").Value = "0.365"
$ws.Range("E10").Value = "  +2.82%  "
$ws.Range("E11").Value = "  -3.83%  "
$ws.Range("D12").Value = "0.0987"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "14.57"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").Value = "2.357.20"
$ws.Range("E14").Value = "  +8.80%  "
$ws.Range("D15").Value = "0.812"
$ws.Range("E15").Value = "  +2.87%  "
$ws.Range("D16").Value = "2.041.29"
$ws.Range("E16").Value = "  +8.23%  "
$ws.Range("D17").Value = "4.92"
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("D18").Value = "36.593.14"
$ws.Range("E18").Value = "  +3.15%  "
$ws.Range("D19").Value = "71.52"
$ws.Range("E19").Value = "  -2.35%  "
$ws.Range("D20").Value = "0.0₃0816"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").Value = "237.99"
$ws.Range("E21").Value = "  -3.31%  "
$ws.Range("D22").Value = "12.65"
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("D23").Value = "4.93"
$ws.Range("E23").Value = "  -4.36%  "
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("E25").Value = "  -7.02%  "
$ws.Range("D26").Value = "169.20"
$ws.Range("E26").Value = "  +2.03%  "
$ws.Range("D27").Value = "20.31"
$ws.Range("E27").Value = "  +10.95%  "
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("E29").Value = "  -8.49%  "
$ws.Range("D30").Value = "0.122"
$ws.Range("E30").Value = "  -4.37%  "
$ws.Range("D31").Value = "21.77"
$ws.Range("E31").Value = "  +50.87%  "
$ws.Range("D32").Value = "4.38"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").Value = "0.0583"
$ws.Range("E33").Value = "  -3.93%  "
$ws.Range("D34").Value = "0.0896"
$ws.Range("E34").Value = "  +19.97%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("D37").Value = "2.27"
$ws.Range("E37").Value = "  +17.21%  "
$ws.Range("D39").Value = "0.871"
$ws.Range("E39").Value = "  +2.52%  "
$ws.Range("E40").Value = "  -10.11%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "1.14"
$ws.Range("E41").Value = "  +4.82%  "
$ws.Range("D42").Value = "96.78"
$ws.Range("E42").Value = "  -1.99%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0215"
$ws.Range("E43").Value = "  -6.01%  "
$ws.Range("E44").Value = "  +15.84%  "
$ws.Range("D45").Value = "16.06"
$ws.Range("E45").Value = "  -5.18%  "
$ws.Range("D46").Value = "1.321.25"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("D47").Value = "0.0816"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("E48").Value = "  +3.11%  "
$ws.Range("D49").Value = "2.244.21"
$ws.Range("E49").Value = "  +8.57%  "
$ws.Range("E50").Value = "  -5.74%  "
$ws.Range("E51").Value = "  +15.69%  "
